$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45916
$ws.Range("B2").Value = 80.01000000000001
$ws.Range("C2").Value = 77.98999999999999
$ws.Range("D2").Value = 72.03
$ws.Range("E2").Value = 73.34999999999999
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 84.2
$ws.Range("H2").Value = 99.90000000000001
$ws.Range("I2").Value = 111.87
$ws.Range("J2").Value = 111.84
$ws.Range("K2").Value = 93.89
$ws.Range("L2").Value = 38.1
$ws.Range("M2").Value = 15.2
$ws.Range("N2").Value = 9.140000000000001
$ws.Range("O2").Value = 9.140000000000001
$ws.Range("P2").Value = 5.76
$ws.Range("Q2").Value = 5.79
$ws.Range("R2").Value = 9.140000000000001
$ws.Range("S2").Value = 38.14
$ws.Range("T2").Value = 69.29000000000001
$ws.Range("U2").Value = 109.03
$ws.Range("V2").Value = 134.29
$ws.Range("W2").Value = 198.95
$ws.Range("X2").Value = 111.84
$ws.Range("Y2").Value = 104.5
$ws.Range("Z2").Value = 72.64
$ws.Range("AB2").Value = 137.39
$ws.Range("AD2").Value = 166.62
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 108.17
$ws.Range("AG2").Value = "2h-18h"

$wb.Save()
